$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.581813931465149
$ws.Range("B1").Value = 2.512386083602905
$ws.Range("C1").Value = 5.672083377838135
$ws.Range("D1").Value = 2.879790544509888
$ws.Range("E1").Value = 0.9197962880134583
